# Pioneer Gliders Calibration and ingest CSV
# Rename "Glider" sheet to "Moorings", update the GL389 glider deployment
# row on the Moorings sheet (launch time, water depth, computed decimal
# lat/long), and repoint the Asset_Cal_Info Ref Des rows from the GL001
# template asset IDs to the real GL389 asset IDs.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Rename the "Glider" sheet to "Moorings"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Glider")
$ws1.Name = "Moorings"

# Renaming the sheet repoints most defined names automatically, but any
# name that referred to a now-invalid range (#REF!) loses its sheet
# qualifier on rename - restore it explicitly.
$wb.Names.Item("_FilterDatabase_0").RefersTo = "=Moorings!#REF!"
$wb.Names.Item("_FilterDatabase_0_0_0").RefersTo = "=Moorings!#REF!"

# ---------------------------------------------------------------------
# 2. Moorings (ex-Glider) sheet: update the GL389 deployment row
# ---------------------------------------------------------------------
$ws1.Range("E2").Value = 0.0625
$ws1.Range("I2").Value = 0

$ws1.Range("L2").Formula = '=((LEFT(G2,(FIND("°",G2,1)-1)))+(MID(G2,(FIND("°",G2,1)+1),(FIND("''",G2,1))-(FIND("°",G2,1)+1))/60))*(IF(RIGHT(G2,1)="N",1,-1))'
$ws1.Range("M2").Formula = '=((LEFT(H2,(FIND("°",H2,1)-1)))+(MID(H2,(FIND("°",H2,1)+1),(FIND("''",H2,1))-(FIND("°",H2,1)+1))/60))*(IF(RIGHT(H2,1)="E",1,-1))'

$ws1.Range("L2:M2").Font.Name = "Calibri"
$ws1.Range("L2:M2").Font.Size = 11
$ws1.Range("L2:M2").Font.ColorIndex = 1
$ws1.Range("L2:M2").HorizontalAlignment = -4108
$ws1.Range("L2:M2").VerticalAlignment = -4108

# ---------------------------------------------------------------------
# 3. Asset_Cal_Info sheet: repoint Ref Des values from GL001 to GL389
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Asset_Cal_Info")

$ws2.Range("A2").Value = "CP05MOAS-GL389-01-ADCPAM000"
$ws2.Range("A7").Value = "CP05MOAS-GL389-02-FLORTM000"
$ws2.Range("A12").Value = "CP05MOAS-GL389-03-CTDGVM000"
$ws2.Range("A14").Value = "CP05MOAS-GL389-04-DOSTAM000"
$ws2.Range("A16").Value = "CP05MOAS-GL389-05-PARADM000"
$ws2.Range("A18").Value = "CP05MOAS-GL389-00-ENG000000"

# ---------------------------------------------------------------------
# 4. Selections: Asset_Cal_Info remembers A35, but Moorings (sheet1)
#    stays the active/selected tab - so select it last.
# ---------------------------------------------------------------------
$ws2.Range("A35").Select()
$ws1.Range("D15").Select()
